$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.014.75'
$ws.Range("E2").Value = '  -3.02%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.820.75'
$ws.Range("E3").Value = '  -2.10%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.83%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.91'
$ws.Range("E5").Value = '  -2.71%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  -0.68%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4221'
$ws.Range("E7").Value = '  -1.89%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3668'
$ws.Range("E8").Value = '  -2.26%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07216'
$ws.Range("E9").Value = '  -1.90%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8399'
$ws.Range("E10").Value = '  -4.38%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.69'
$ws.Range("E11").Value = '  -4.43%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.823.65'
$ws.Range("E12").Value = '  -1.95%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.663'
$ws.Range("E13").Value = '  -1.44%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07062'
$ws.Range("E14").Value = '  -1.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.282'
$ws.Range("E15").Value = '  -3.23%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '89.98'
$ws.Range("E16").Value = '  +1.06%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.003'
$ws.Range("E17").Value = '  -0.97%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008730'
$ws.Range("E18").Value = '  -3.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  -0.79%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.88'
$ws.Range("E20").Value = '  -3.91%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.168.58'
$ws.Range("E21").Value = '  -2.46%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.140'
$ws.Range("E22").Value = '  -1.59%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.81'
$ws.Range("E23").Value = '  -2.73%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.056.11'
$ws.Range("E24").Value = '  -1.21%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.986'
$ws.Range("E25").Value = '  +0.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.76'
$ws.Range("E26").Value = '  -2.34%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.269'
$ws.Range("E27").Value = '  +4.18%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.28'
$ws.Range("E28").Value = '  -2.07%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.255'
$ws.Range("E29").Value = '  -2.38%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '116.87'
$ws.Range("E30").Value = '  -1.76%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08705'
$ws.Range("E31").Value = '  -2.74%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.173'
$ws.Range("E32").Value = '  -5.12%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7330'
$ws.Range("E33").Value = '  -6.17%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.889'
$ws.Range("E34").Value = '  -1.28%  '
$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.410'
$ws.Range("E35").Value = '  -3.22%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.001'
$ws.Range("E36").Value = '  -0.85%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.090'
$ws.Range("E37").Value = '  -3.57%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01949'
$ws.Range("E38").Value = '  -1.87%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05250'
$ws.Range("E39").Value = '  -2.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.297'
$ws.Range("E40").Value = '  +0.78%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.882'
$ws.Range("E41").Value = '  -0.95%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1683'
$ws.Range("E42").Value = '  -1.11%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.5051'
$ws.Range("E43").Value = '  -2.00%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.543'
$ws.Range("E44").Value = '  -3.80%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.52'
$ws.Range("E45").Value = '  -2.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '105.98'
$ws.Range("E46").Value = '  -2.26%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4703'
$ws.Range("E47").Value = '  -1.54%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.926'
$ws.Range("E48").Value = '  +3.71%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.001'
$ws.Range("E49").Value = '  -0.86%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06317'
$ws.Range("E50").Value = '  -2.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.647'
$ws.Range("E51").Value = '  -2.91%  '
